$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAIO")

# Insert a new row above the current row 6 (shifts existing rows 6.. down by one),
# copying the formatting of the row above (row 5) so the new row matches its
# neighbours instead of picking up a fresh style.
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(6).Insert()

$ws.Cells.Item(6, 1).Value = "JHSF STATES - ERICA - TRATOR DE ESTEIRA FOI DIA 29/05"
$ws.Cells.Item(6, 2).Formula = "=(35000/31)*3"
$ws.Cells.Item(6, 3).Value = $null
